$d = $word.ActiveDocument

$replacements = @(
    @("2024-02-04 Sunday", "2024-02-05 Monday"),
    @("90×47=4230", "52×15=780"),
    @("83×85=7055", "94×11=1034"),
    @("14×28=392", "99×39=3861"),
    @("25×64=1600", "37×32=1184"),
    @("19×85=1615", "45×76=3420"),
    @("19×79=1501", "74×78=5772"),
    @("66×95=6270", "69×70=4830"),
    @("24×69=1656", "64×25=1600"),
    @("47×81=3807", "30×96=2880"),
    @("99×52=5148", "59×45=2655"),
    @("20×36=720", "84×52=4368"),
    @("28×81=2268", "96×87=8352"),
    @("49×75=3675", "70×63=4410"),
    @("88×92=8096", "62×79=4898"),
    @("12×80=960", "88×78=6864"),
    @("38×31=1178", "28×85=2380"),
    @("93×24=2232", "88×83=7304"),
    @("56×85=4760", "49×97=4753"),
    @("77×67=5159", "90×97=8730"),
    @("65×78=5070", "81×69=5589"),
    @("64×87=5568", "14×63=882"),
    @("68×91=6188", "71×32=2272"),
    @("95×49=4655", "16×28=448"),
    @("65×14=910", "52×31=1612"),
    @("58×72=4176", "64×54=3456")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
